$wb = $excel.ActiveWorkbook

# The "Coding" sheet (tab 2) is the one being edited - it is already the
# active sheet/tab in the workbook (tabSelected="1", activeTab index 1).
$ws = $wb.Worksheets.Item("Coding")

# Remove the "without authors" note from the book-form row (row 4) and the
# stray "x"/"without authors" marks that no longer apply now that the
# "NewBookWindow" supports an authors selection - clear the cell contents
# but keep the existing cell formatting (borders) in place.
$ws.Range("E4").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("E7").ClearContents()

# Update the active selection on the sheet to reflect where the user ended
# up after making the edit.
$ws.Range("G5").Select()
